$wb = $excel.ActiveWorkbook

# --- Netherlands: copy of "Greece" (closest template), placed at the end ---
$greece = $wb.Worksheets.Item("Greece")
$greece.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$netherlands = $wb.Worksheets.Item($wb.Worksheets.Count)
$netherlands.Name = "Netherlands"
$netherlands.Range("B4").Value = "NGC-3144/T2176"
$netherlands.Range("B2").Value = "Netherlands Market"
$netherlands.Rows.Item(2).RowHeight = 28.8
$netherlands.Rows.Item(3).RowHeight = 28.8
$netherlands.Rows.Item(4).RowHeight = 28.8
$netherlands.Range("B4").Select() | Out-Null

# --- Austria: copy of "Croatia" (exact structural template match) ---
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$austria = $wb.Worksheets.Item($wb.Worksheets.Count)
$austria.Name = "Austria"
$austria.Range("B4").Value = "NGC-3817/T2272"
$austria.Range("B2").Value = "Austria Market"
$austria.Range("B4").Select() | Out-Null

# --- Denmark: copy of "Croatia" (exact structural template match) ---
$croatia.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B4").Value = "NGC-2913/T2749"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Select() | Out-Null

# Denmark becomes the active/selected tab, matching the new activeTab in the workbook view.
$denmark.Activate() | Out-Null
